$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last data row (A6) onto the new row so the
# new cell reuses the existing "s=2" style instead of creating a new one
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new city entry: PS / Pachmarhi
$ws.Range("A7").Value = "PS"
$ws.Range("B7").Value = "Pachmarhi"

# Move the active selection down to the next empty row, as in the saved file
$ws.Range("B8").Select()
